$d = $word.ActiveDocument

# Locate the paragraph that contains the "Ver no Jupiter..." text and the
# paragraph that contains the "(c) 2020 ..." text, plus the blank paragraph
# that sits between "LOM3096: ..." and "Ver no Jupiter ...". All three
# paragraphs are removed as a single contiguous block.

$startPara = $null
$endPara = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $startPara = $p
    }
    if ($t -like "*Contact: luizeleno@usp.br*") {
        $endPara = $p
    }
}

# Expand the start to include the blank paragraph right before "Ver no Jupiter ..."
$blockStart = $startPara.Previous().Range.Start
$blockEnd = $endPara.Range.End

$r = $d.Range($blockStart, $blockEnd)
$r.Delete()
